$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns D..G ---
$ws.Range("D1").Value = "FSKONZ"
$ws.Range("E1").Value = "FSFIRM"
$ws.Range("F1").Value = "FSAN1"
$ws.Range("G1").Value = "BESTANDSMODULE"

# Copy the existing header style (bold font + grey fill) from C1 onto the
# newly added header cells E1:G1 (D1 already carries it from the source file).
$ws.Range("C1").Copy() | Out-Null
$ws.Range("E1:G1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Data rows 2-14 ---
$data = @(
    @("HV9","TST","E09","D01","040","                                   ",34),
    @("HV9","TST","E09","S01","020","Google Germany GmbH                ",9),
    @("HV9","TST","E09","D01","030","Test F.Bäcker                      ",8686),
    @("HV9","TST","E09","D12","010","Test AKL                           ",95),
    @("HV9","TST","E09","D97","010","FIEGE Logistik Stiftung & Co. KG   ",9),
    @("HV9","TST","E09","S01","SME","Sascha Mergard                     ",77),
    @("HV9","TST","E09","D14","010","JBE                                ",107),
    @("HV9","TST","E09","F01","010","Firat Askin                        ",2),
    @("HV9","TST","E09","D80","FFA","A                                  ",66),
    @("HV9","TST","E09","D02","010",".                                  ",146),
    @("HV9","TST","E09","D13","BWO","Bastian Woltemade GmbH             ",57),
    @("HV9","TST","E09","D99","000","Schulung                           ",7),
    @("HV9","TST","E09","D03","100","KST TESTFIRMA                      ",148)
)

# Column E holds codes such as "040" / "020" / "010" that must stay text
# (leading zeros), so force text format before writing the values, then put
# the cell style back to Normal so the saved file doesn't carry a stray
# "text number format" style on these cells (matches the source workbook,
# which has no explicit style on the data rows).
$ws.Range("E2:E14").NumberFormat = "@"

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $row = $row + 1
}

$ws.Range("E2:E14").Style = "Normal"

# --- Column widths for the new columns D..G ---
# Empirically, this host's ColumnWidth -> stored <col width="..."> conversion
# adds 11/12 on write, so subtract that to land exactly on the target widths
# (8, 8, 37, 16) in the saved OOXML.
$ws.Range("D1").ColumnWidth = 8 - 11/12
$ws.Range("E1").ColumnWidth = 8 - 11/12
$ws.Range("F1").ColumnWidth = 37 - 11/12
$ws.Range("G1").ColumnWidth = 16 - 11/12

# --- AutoFilter: extend from A1:D1 to A1:G1 ---
# The sheet already has AutoFilterMode on (ref A1:D1); toggle it off first so
# re-applying AutoFilter() on the wider range turns it back on with the new ref.
$ws.AutoFilterMode = $False
$ws.Range("A1:G1").AutoFilter() | Out-Null

# --- Update the _xlnm._FilterDatabase defined name to match the new range ---
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='SQL Ergebnisse'!`$A`$1:`$G`$1"
    }
}
